$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 1).Value = "تب کیو توسط کدام عامل ایجاد می‌شود؟"
$ws.Cells.Item(3, 1).Value = "در سندرم مارفان، آنوریسم آئورت در کدام قسمت رخ می‌دهد؟"
$ws.Cells.Item(4, 1).Value = "رنگ گرانول اکتینومیست‌ها چیست؟"
$ws.Cells.Item(5, 1).Value = "دیسپلازی بال اسفنوئید در کدام یک از موارد زیر دیده می‌شود؟"
$ws.Cells.Item(6, 1).Value = "مرکز وازوموتور مدولا با کدام یک مرتبط است؟"
$ws.Cells.Item(7, 1).Value = "بیمار دچار درماتیت حساس به نور خورشید، اسهال و زوال عقل است. کمبود کدام یک از موارد زیر باید در او وجود داشته باشد؟"
$ws.Cells.Item(8, 1).Value = "بیماری موندور شامل کدام یک از موارد زیر می‌شود؟"
$ws.Cells.Item(9, 1).Value = "جایگزینی اپیتلیوم ستونی دستگاه تنفسی با اپیتلیوم سنگفرشی چیست؟"
$ws.Cells.Item(10, 1).Value = "فعال کننده کامپوزیت نوری چیست؟"
$ws.Cells.Item(11, 1).Value = "توموری که پسرفت خودبه‌خودی نشان نمی‌دهد:"
$ws.Cells.Item(12, 1).Value = "استخوان های سوار (rider's bones) استخوانی شدن کدام عضله هستند؟"
$ws.Cells.Item(13, 1).Value = "در یک مورد سرطان دهانه رحم، بیمار با تغییر سطح هوشیاری و سکسکه مراجعه می‌کند. علت احتمالی چیست؟"
$ws.Cells.Item(14, 1).Value = "مدیریت پارگی داخل صفاقی مثانه:"
$ws.Cells.Item(15, 1).Value = "در مورد سندرم فمینیزاسیون بیضه کدام گزینه صحیح است؟"
$ws.Cells.Item(16, 1).Value = "یک مرد 38 ساله از کاهش شنوایی در گوش راست در طی دو سال گذشته شکایت دارد. در آزمایش با چنگال تنظیم 512 هرتز، تست رینه بدون ماسک در گوش راست منفی و در گوش چپ مثبت است. در تست وبر، صدا در گوش چپ بلندتر شنیده می‌شود. این بیمار به احتمال زیاد چه مشکلی دارد؟"
$ws.Cells.Item(17, 1).Value = "در مورد شکستگی ترقوه کدام گزینه صحیح است؟  
الف) جوش خوردن بدشکل  
ب) شایع‌ترین محل شکستگی یک‌سوم داخلی و دو‌سوم خارجی است  
ج) شکستگی خرد شده  
د) ناشی از افتادن روی دست کشیده‌شده"
$ws.Cells.Item(18, 1).Value = "هتروکرومیای عنبیه از ویژگی‌های کدام یک از موارد زیر است؟"
$ws.Cells.Item(19, 1).Value = "کدام یک از تاندون‌های زیر از زیر سوستانتاکولوم تالی عبور می‌کند؟"
$ws.Cells.Item(20, 1).Value = "کمبود آنزیم b-glucosidase باعث کدام بیماری می‌شود؟"
$ws.Cells.Item(21, 1).Value = "روش ترانس پریتوانیال برای ترمیم فیستول ژنیتوریناری در تمام موارد زیر به جز کدام مورد استفاده می‌شود؟"
$ws.Cells.Item(22, 1).Value = "کدام یک از این آنزیم های گوارشی در شیره لوزالمعده وجود ندارد؟"
$ws.Cells.Item(23, 1).Value = "سلول‌های 'کلو' در کدام مورد مشاهده می‌شوند؟"
$ws.Cells.Item(24, 1).Value = "ادم غیر گوده‌ای در کدام مورد مشاهده می‌شود؟"
$ws.Cells.Item(25, 1).Value = "یک نوزاد تازه متولد شده در بخش نوزادان دچار افسردگی تنفسی شده است. کدام یک از داروهای زیر علت این وضعیت است؟"
$ws.Cells.Item(26, 1).Value = "عدم تشکیل ماندیبل و ناشنوایی هدایتی ممکن است به دلیل عدم تشکیل کدام قوس حلقی باشد؟"
$ws.Cells.Item(27, 1).Value = "در یک کودک پنج ساله، کدام یک از سوراخ‌های زیر در حفره جمجمه‌ای خلفی قرار دارد؟"
$ws.Cells.Item(28, 1).Value = "جراحی در کدام مورد مفید نیست؟"
$ws.Cells.Item(29, 1).Value = "شایع ترین دفورمیتی مشاهده شده در استئوآرتریت اولیه مفصل زانو کدام است؟"
$ws.Cells.Item(30, 1).Value = "همه موارد زیر در درمان آبسه کبدی آمیبی استفاده می‌شوند به جز -"
$ws.Cells.Item(31, 1).Value = "کدام یک از موارد زیر در مورد فتق اپی گاستریک صحیح است؟"
$ws.Cells.Item(32, 1).Value = "گلیکوژنین پروتئینی با قابلیت خودگلیکوزیلاسیون است. مولکول‌های گلوکز به کدام آمینواسید گلیکوژنین متصل می‌شوند؟"
$ws.Cells.Item(33, 1).Value = "در مورد تورهای پیشگیری از پشه، همه موارد زیر صحیح هستند به جز"
$ws.Cells.Item(34, 1).Value = "یک زن ۵۰ ساله با توده بافت نرم بدون درد در ران راست خود مراجعه می‌کند. پس از برداشتن توده با جراحی، جراح متوجه چسبندگی تومور به بافت‌های اطراف می‌شود. بررسی هیستولوژیک نشان‌دهنده نئوپلاسمی متشکل از سلول‌های پلئومورفیک با سیتوپلاسم واکوئله و شفاف است. هسته بسیاری از سلول‌ها توسط واکوئل‌های سیتوپلاسمی که با روش‌های هیستوشیمیایی برای لیپیدها رنگ می‌گیرند، فرورفته است. کدام یک از موارد زیر محتمل‌ترین تشخیص است؟"
$ws.Cells.Item(35, 1).Value = "نقص در اراده معمولاً در کدام یک از موارد زیر مشاهده می‌شود؟"
$ws.Cells.Item(36, 1).Value = "سنتز کدام یک از فاکتورهای انعقادی زیر در بیماری کبدی تحت تأثیر قرار نمی‌گیرد؟"
$ws.Cells.Item(37, 1).Value = "در نرخ مرگ و میر نوزادان (imr)، سن نوزادان در نظر گرفته می‌شود –"
$ws.Cells.Item(38, 1).Value = "کدام یک از گزینه‌های زیر دقیق‌ترین توصیف از تونسیلکتومی داخل کپسولی کوبلیشن را ارائه می‌دهد؟"
$ws.Cells.Item(39, 1).Value = "اجسام اسکلروتیک در کدام یک از موارد زیر دیده می‌شود؟"
$ws.Cells.Item(40, 1).Value = "اسکالوگرام نام دیگر کدام مورد است؟"
$ws.Cells.Item(41, 1).Value = "در درمان اندودنتیک بازساختی دندان‌های دائمی جوان، مترونیدازول با کدام دارو جایگزین می‌شود؟"
$ws.Cells.Item(42, 1).Value = "کدام یک از موارد زیر در مورد مرحله زیر از احیای نوزادان صحیح است؟"
$ws.Cells.Item(43, 1).Value = "مانیتول برای کدام مورد مفید نیست؟"
$ws.Cells.Item(44, 1).Value = "دنتین بین گلوبولی در نتیجه کدام مورد ایجاد می‌شود؟"
$ws.Cells.Item(45, 1).Value = "کدام یک از موارد زیر بر اساس معیارهای who در سال 2010 برای آنالیز طبیعی مایع منی صحیح نیست؟"
$ws.Cells.Item(46, 1).Value = "مری حلزونی (تسبیحی) به طور مشخص در کدام یک از موارد زیر دیده می‌شود؟"
$ws.Cells.Item(47, 1).Value = "عاملی که با افزودن به یک کلونی، رشد آن را مهار میکند و با حذف آن، کلونی مجدداً رشد میکند، کدام است؟"
$ws.Cells.Item(48, 1).Value = "کدام یک از موارد زیر درماتوم ناف است؟"
$ws.Cells.Item(49, 1).Value = "هنگامی که سلول‌های بنیادی به سلول‌های مشخصه بافت‌های دیگر تبدیل می‌شوند، این فرآیند چه نامیده می‌شود؟"
$ws.Cells.Item(50, 1).Value = "غدد لنفاوی که در سرطان دهانه رحم درگیر نمی‌شوند"
$ws.Cells.Item(51, 1).Value = "ویژگی‌های بالینی دندان‌های شیری همه موارد زیر هستند به جز:"
$ws.Cells.Item(52, 1).Value = "دمای تخته مخلوط کردن برای سیمان های سیلیکاتی (بر حسب فارنهایت) باید چقدر باشد؟"
$ws.Cells.Item(53, 1).Value = "مرد 25 ساله با آسیب به سر که در حال راه رفتن است، هوشیار و جهت‌یافته می‌باشد و سپس بیهوش می‌شود. تشخیص چیست؟"
$ws.Cells.Item(54, 1).Value = "پرتو درمانی با دوز پایین در مدت زمان کمتر باعث همه موارد زیر می‌شود به جز؟"
$ws.Cells.Item(55, 1).Value = "آرژینین با تشکیل کدام سوبسترا وارد چرخه tca می‌شود؟"
$ws.Cells.Item(56, 1).Value = "در یک جمعیت 5000 نفری، 500 نفر از یک بیماری رنج می‌برند. همچنین 150 مورد قدیمی از این بیماری در جمعیت وجود دارد. میزان شیوع بیماری را محاسبه کنید."
$ws.Cells.Item(57, 1).Value = "شایع‌ترین علت هموبیلیا چیست؟"
$ws.Cells.Item(58, 1).Value = "یک بیمار با هایپرتریگلیسریدمی با اسیدهای چرب چند غیراشباع امگا-3 درمان می‌شود. درمان با اسیدهای چرب چند غیراشباع امگا-3، چه تأثیری بر پروفایل لیپیدی خواهد داشت؟"
$ws.Cells.Item(59, 1).Value = "یک مرد 40 ساله با سابقه مصرف روزانه الکل در طول 7 سال گذشته، با شروع حاد دیدن مارها در اطراف خود در اتاق، عدم شناسایی اعضای خانواده، رفتار خشونت‌آمیز و لرزش برای چند ساعت به اورژانس بیمارستان آورده شده است. سابقه عدم مصرف الکل در 2 روز گذشته وجود دارد. معاینه نشان‌دهنده افزایش فشار خون، لرزش، افزایش فعالیت روانی-حرکتی، تأثیر ترس، رفتار توهمی، گم‌گشتگی، اختلال در قضاوت و بینش است. او به احتمال زیاد از چه بیماری رنج می‌برد؟"
$ws.Cells.Item(60, 1).Value = "کدام گزینه در مورد سندرم چرگ اشتراوس نادرست است؟"
$ws.Cells.Item(61, 1).Value = "درد سوزشی اپی گاستر ناشی از چیست؟"
$ws.Cells.Item(62, 1).Value = "شوک در بالین بهتر است توسط کدام مورد ارزیابی شود؟"
$ws.Cells.Item(63, 1).Value = "کیست بیکر چیست؟"
$ws.Cells.Item(64, 1).Value = "درمان ناباروری سرویکال می‌تواند شامل همه موارد زیر باشد به جز:"
$ws.Cells.Item(65, 1).Value = "همه موارد زیر از خانواده فلاویویروس‌ها هستند به جز،"
$ws.Cells.Item(66, 1).Value = "گسترش لثه برای یک ترمیم باید در کجا باشد؟"
$ws.Cells.Item(67, 1).Value = "تغییرات در خون هنگام عبور از مویرگ‌های سیستمیک شامل همه موارد زیر است به جز"
$ws.Cells.Item(68, 1).Value = "کدام یک از موارد زیر با بیماری کلیه پلی‌کیستیک بزرگسالان مرتبط است؟"
$ws.Cells.Item(69, 1).Value = "همه موارد زیر در مورد تبلیغات سلامت صحیح هستند به جز"
$ws.Cells.Item(70, 1).Value = "یک زن باردار مبتلا به بیماری گریوز تشخیص داده شده است. مناسب‌ترین درمان برای او کدام است؟"
$ws.Cells.Item(71, 1).Value = "یافته های میلوم مولتیپل در کلیه شامل همه موارد زیر به جز کدام است؟"
$ws.Cells.Item(72, 1).Value = "لیپاز حساس به هورمون توسط کدام مورد فعال نمی‌شود؟"
$ws.Cells.Item(73, 1).Value = "کمبود g-6-pd به چه صورت به ارث می رسد؟"
$ws.Cells.Item(74, 1).Value = "آسیب ناشی از وارونگی اجباری در پای در حالت پلانتار فلکشن کدام است؟"
$ws.Cells.Item(75, 1).Value = "طوفان تیروئیدی پس از عمل به دلیل چیست؟"
$ws.Cells.Item(76, 1).Value = "تحریک واگ باعث اثرات زیر می‌شود به جز"
$ws.Cells.Item(77, 1).Value = "در بتا تالاسمی، شایع ترین جهش ژنی کدام است؟"
$ws.Cells.Item(78, 1).Value = "زولپیدم-"
$ws.Cells.Item(79, 1).Value = "در طبقه‌بندی فارست برای زخم معده خونریزی‌دهنده با رگ قابل مشاهده یا برجستگی رنگدانه‌دار، کدام گزینه صحیح است؟"
$ws.Cells.Item(80, 1).Value = "درباره داروی اسیدی کدام گزینه صحیح است -"
$ws.Cells.Item(81, 1).Value = "ماده مهارکننده مولرین توسط کدام یک ترشح می‌شود؟"
$ws.Cells.Item(82, 1).Value = "همه موارد زیر از ویژگی‌های بالینی بیماری التهابی لگن (pid) هستند به جز:"
$ws.Cells.Item(83, 1).Value = "بیمار با کارسینوم سلول کلیوی که به ivc و ورید کلیوی تهاجم کرده است، مراجعه کرده است. گزینه نادرست کدام است؟"
$ws.Cells.Item(84, 1).Value = "تست لپرومین در چه زمانی مشاهده می‌شود؟"
$ws.Cells.Item(85, 1).Value = "کدام یک از موارد زیر جزو تست‌های قطعی برای پایش رشد و تکامل محسوب نمی‌شود؟"
$ws.Cells.Item(86, 1).Value = "پدیده پروزون به دلیل چیست؟"
$ws.Cells.Item(87, 1).Value = "بیماری تای-ساکس به دلیل کمبود کدام یک از موارد زیر رخ می‌دهد؟"
$ws.Cells.Item(88, 1).Value = "در عفونت کزاز:"
$ws.Cells.Item(89, 1).Value = "یک نوزاد پسر 3 ماهه دچار اوتیت مدیا شد که برای آن یک دوره کو-تریموکسازول دریافت کرد. چند روز بعد، پوست او به طور گسترده شروع به پوست انداختن کرد؛ هیچ ضایعه مخاطی وجود نداشت و نوزاد حال عمومی بدی نداشت. محتمل‌ترین تشخیص چیست؟"
$ws.Cells.Item(90, 1).Value = "کدام یک از موارد زیر باید در ابتدا به یک کودک دچار سوءتغذیه شدید داده شود؟"
$ws.Cells.Item(91, 1).Value = "در سندرم‌های men نوع 2a و 2b، کدام گیرنده فاکتور رشد دچار جهش نقطه‌ای می‌شود؟"
$ws.Cells.Item(92, 1).Value = "فعال شدن کدام ژن منجر به کارسینوم مثانه می‌شود؟"
$ws.Cells.Item(93, 1).Value = "کدام مورد در نمره زخم سپسیس شامل نمی‌شود؟"
$ws.Cells.Item(94, 1).Value = "اسفنکتر مجرای ادرار در کدام قسمت وجود دارد؟"
$ws.Cells.Item(95, 1).Value = "یافته غلظت سدیم ادرار کمتر از 20 میلی‌مول بر لیتر در هیپوناترمی هیپوولمیک نشان‌دهنده کدام مورد است؟"
$ws.Cells.Item(96, 1).Value = "کدام گزینه در مورد انتقال پلاکت صحیح نیست؟"
$ws.Cells.Item(97, 1).Value = "یک دانشجوی ۲۳ ساله پس از پریدن از یک آبشار ۵۰ فوتی به بخش اورژانس منتقل می‌شود. تصویربرداری mri از کمر او نشان‌دهنده جابجایی جانبی نخاع به سمت چپ است. کدام یک از ساختارهای زیر به احتمال زیاد پاره شده‌اند تا باعث این انحراف شوند؟"
$ws.Cells.Item(98, 1).Value = "`"اجسام اگزنر`" در کدام مورد دیده می‌شوند؟"
$ws.Cells.Item(99, 1).Value = "دررفتگی شکستگی بنت شامل پایه کدام است؟"
$ws.Cells.Item(100, 1).Value = "کدام ناحیه خارج از سد خونی مغزی قرار دارد؟"
$ws.Cells.Item(101, 1).Value = "در صورت سوختگی با ضخامت جزئی، چه زمانی بیمار را به مرکز سوختگی ارجاع می‌دهید؟"
$ws.Cells.Item(102, 1).Value = "داروی مورد استفاده برای واکنش دیستونیک ناشی از متوکلوپرامید کدام است؟"
$ws.Cells.Item(103, 1).Value = "مجرای پانکراس شکمی به کدام بخش منتهی می‌شود؟"
$ws.Cells.Item(104, 1).Value = "کدام یک از عوامل زیر باعث ایجاد اولسر سرپنس مشخص قرنیه می شود؟"
$ws.Cells.Item(105, 1).Value = "شایع‌ترین عمل جراحی انجام‌شده برای زخم معده با انسداد خروجی معده چیست؟"
$ws.Cells.Item(106, 1).Value = "نادرست در مورد سندرم روده تحریک پذیر"
$ws.Cells.Item(107, 1).Value = "دوره ایمن کدام است؟"
$ws.Cells.Item(108, 1).Value = "یک مرد مسن با همی‌پارزی آتاکسیک، تشخیص چیست؟"
$ws.Cells.Item(109, 1).Value = "حلقه ویمبرگر در کدام بیماری دیده می‌شود؟ سپتامبر 2010"
$ws.Cells.Item(110, 1).Value = "ناتوانی جنسی می‌تواند به عنوان دفاع در مردان در تمام شرایط زیر مطرح شود به جز"
$ws.Cells.Item(111, 1).Value = "دیسمنوره احتقانی در بیماران مبتلا به کدام یک از موارد زیر مشاهده می‌شود؟"
$ws.Cells.Item(112, 1).Value = "شایع‌ترین محل بروز کارسینوم سلول سنگفرشی در حفره دهان کدام است؟"
$ws.Cells.Item(113, 1).Value = "عملکرد مخچه مغزی (cerebrocerebellum) چیست؟"
$ws.Cells.Item(114, 1).Value = "'لرزش های هتر' در مسمومیت با ... دیده می‌شود:"
$ws.Cells.Item(115, 1).Value = "سندرم زیر با برونشکتازی همراه است -"
$ws.Cells.Item(116, 1).Value = "مرگ مغزی زمانی رخ می‌دهد که: مارس 2008"
$ws.Cells.Item(117, 1).Value = "کدام یک از موارد زیر نمونه‌ای از مقیاس اسمی نیست؟"
$ws.Cells.Item(118, 1).Value = "کشیدن دندان مولر اول فک پایین در یک کودک 8 ساله به عنوان چه چیزی شناخته می‌شود؟"
$ws.Cells.Item(119, 1).Value = "شریان کاروتید مشترک در کدام مرز غضروف تیروئید دو شاخه می‌شود؟"
$ws.Cells.Item(120, 1).Value = "شایع ترین عامل ایجاد کننده استئومیلیت در کم خونی داسی شکل کدام است؟"
$ws.Cells.Item(121, 1).Value = "مفصل میانی آتلانتو-آکسیال از چه نوع مفصلی است؟"
$ws.Cells.Item(122, 1).Value = "کدام یک از دوزهای زیر در واحدهای لوفلر توکسوئید دیفتری در هر دوز واکسن dpt گنجانده شده است؟"
$ws.Cells.Item(123, 1).Value = "علت لنفادنیت مزانتریک حاد چیست؟"
$ws.Cells.Item(124, 1).Value = "در مورد مایکوباکتریوم غیر از توبرکلوزیس کدام گزینه صحیح است؟"
$ws.Cells.Item(125, 1).Value = "کدام عضله جزو دیافراگم لگنی محسوب نمی‌شود؟"
$ws.Cells.Item(126, 1).Value = "روش تشخیصی انتخابی در کله دوکولیتیازیس چیست؟"
$ws.Cells.Item(127, 1).Value = "نام روش نشان داده شده در حالب را ذکر کنید:"
$ws.Cells.Item(128, 1).Value = "ناقلین ارگانیسم های غیر بیماری زا به عنوان چه شناخته می شوند؟"
$ws.Cells.Item(129, 1).Value = "مهم‌ترین عامل خطر برای ایجاد عفونت/سپسیس رحم پس از زایمان چیست؟"
$ws.Cells.Item(130, 1).Value = "یک بیمار با نقص دیواره بین بطنی (vsd) دچار فشار خون ریوی می‌شود، ویژگی مشخصه آن چیست؟"
$ws.Cells.Item(131, 1).Value = "به طور کلی، تیز کردن دستی ابزارها با سنگ‌های بدون پایه نسبت به تیز کردن با سنگ‌های موتوردار ترجیح داده می‌شود زیرا سنگ‌های بدون پایه"
$ws.Cells.Item(132, 1).Value = "در انسداد مجاری هوایی فوقانی تمامی تغییرات زیر مشاهده می‌شود به جز:"
$ws.Cells.Item(133, 1).Value = "آسیت موکینی در کدام مورد مشاهده می‌شود؟"
$ws.Cells.Item(134, 1).Value = "تمام داروهای زیر در بیماران مبتلا به فئوکروموسیتوم باعث بحران فشار خون می‌شوند، به جز:"
$ws.Cells.Item(135, 1).Value = "ضخیم‌ترین بخش صلبیه (اسکلرا) کدام است؟"
$ws.Cells.Item(136, 1).Value = "درمان قطعی در گلوکوم با زاویه بسته اولیه: مارس 2005"
$ws.Cells.Item(137, 1).Value = "تمام رفلکس‌های زیر در بدو تولد وجود دارند به جز _______"
$ws.Cells.Item(138, 1).Value = "آقای سونیل، مرد ۲۲ ساله‌ای که پس از درگیری در خیابان توسط یک فرد مست به صورتش مشت زده شده بود، به اورژانس آورده شد. بیمار هوشیار بود و از گوش راستش خون جاری بود. خونریزی در آن زمان کنترل شد. روز بعد، آزمون شنوایی‌سنجی بیمار نشان‌دهنده کاهش شنوایی از نوع انتقالی بود. نوع شکستگی استخوان گیجگاهی که بیمار دارد چیست؟"
$ws.Cells.Item(139, 1).Value = "آزمایش‌های کنترل‌شده تصادفی‌شده عبارتند از -"
$ws.Cells.Item(140, 1).Value = "کم خونی ماکروسیتیک توسط همه موارد زیر به جز ایجاد می شود:"
$ws.Cells.Item(141, 1).Value = "ریشه‌های درگیر در فلج ارب کدامند؟"
$ws.Cells.Item(142, 1).Value = "روش تشخیصی انتخابی برای یک زن باردار با توده شکمی فوقانی چیست؟"
$ws.Cells.Item(143, 1).Value = "اسیدهای چرب توسط همه موارد زیر به جز کدام مورد استفاده می‌شوند؟"
$ws.Cells.Item(144, 1).Value = "در تصویربرداری اشعه ایکس معمولی از شانه، کدام ساختار در بالاترین موقعیت قرار دارد؟"
$ws.Cells.Item(145, 1).Value = "در مورد سندرم چدیاک-هیگاشی، کدام گزینه صحیح است؟"
$ws.Cells.Item(146, 1).Value = "عمق آب‌بندی در توالت rca چقدر است؟"
$ws.Cells.Item(147, 1).Value = "در هفته ۲۵ بارداری، یک زن g1p0 دچار افزایش فشار خون (۱۶۰/۹۵) و پروتئینوری می‌شود. معاینه فیزیکی ادم عمومی را نشان می‌دهد و آزمایشات سرمی هایپراوریسمی و افزایش غلظت آنزیم‌های کبدی را نشان می‌دهند. کدام یک از موارد زیر معمولاً درمان قطعی برای این وضعیت پزشکی بیمار ارائه می‌دهد؟"
$ws.Cells.Item(148, 1).Value = "عصب‌دهی پرده صماخ توسط کدام است؟"
$ws.Cells.Item(149, 1).Value = "lathi می‌تواند باعث تمام آسیب‌های زیر شود، به جز:"
$ws.Cells.Item(150, 1).Value = "درد شدید مرتبط با عفونت انگشت (felon) به چه دلیل رخ می‌دهد؟"
$ws.Cells.Item(151, 1).Value = "شایع‌ترین محل کارسینوم نازوفارنکس کدام است؟"
$ws.Cells.Item(152, 1).Value = "کدام یک از موارد زیر یک داروی ضد سرطان غیر وابسته به چرخه سلولی است؟"
$ws.Cells.Item(153, 1).Value = "در مورد baha کدام گزینه صحیح است؟"
$ws.Cells.Item(154, 1).Value = "اندام های لنفاوی مرکزی کدامند؟"
$ws.Cells.Item(155, 1).Value = "افزایش سطح 5-hiaa در کدام بیماری مشاهده می‌شود؟"
$ws.Cells.Item(156, 1).Value = "انطباق پذیری ریه ها چقدر است؟"
$ws.Cells.Item(157, 1).Value = "کدام سلول بیشترین حساسیت را به پرتودرمانی دارد؟"
$ws.Cells.Item(158, 1).Value = "شایع‌ترین عارضه شکستگی کولز چیست؟"
$ws.Cells.Item(159, 1).Value = "گلوکونئوژنز در تمام موارد زیر به جز کدام یک رخ می‌دهد؟"
$ws.Cells.Item(160, 1).Value = "آزمون رینه در کدام مورد مثبت است؟"
$ws.Cells.Item(161, 1).Value = "بازداری متقابل عضله آنتاگونیست در نگاه جانبی توسط کدام قانون توضیح داده می‌شود؟"
$ws.Cells.Item(162, 1).Value = "کاردیومیوپاتی پری پارتوم در چه زمانی رخ می‌دهد؟"
$ws.Cells.Item(163, 1).Value = "شکاف لب، پلی‌داکتیلی و هولوپروزنسفالی از ویژگی‌های کدام مورد هستند؟"
$ws.Cells.Item(164, 1).Value = "کوتاه‌ترین طول گوش میانی:"
$ws.Cells.Item(165, 1).Value = "کدام یک از موارد زیر از طریق گیرنده تیروزین کیناز عمل می‌کند؟"
$ws.Cells.Item(166, 1).Value = "شایع‌ترین عارضه شکستگی کولز چیست؟"
$ws.Cells.Item(167, 1).Value = "آب مروارید مادرزادی یک طرفه کامل ترجیحاً باید در چه زمانی عمل شود؟"
$ws.Cells.Item(168, 1).Value = "نسبت آب به پودر برای سنگ دندانی و گچ به ترتیب چیست؟"
$ws.Cells.Item(169, 1).Value = "میکروسکوپی زمینه تاریک برای تشخیص کدام مورد استفاده می‌شود؟"
$ws.Cells.Item(170, 1).Value = "یک زن ۴۵ ساله با درد گردن در بیمارستان بستری می‌شود. سی‌تی‌اسکن یک تومور در سمت چپ حفره دهان او را نشان می‌دهد. تومور و بافت‌های مرتبط با آن با یک عمل جراحی رادیکال گردن برداشته می‌شوند. دو ماه پس از عمل، شانه چپ بیمار به‌طور محسوسی افتاده است. معاینه فیزیکی ضعف واضح در چرخاندن سر به سمت راست و اختلال در abduction اندام فوقانی چپ تا سطح شانه را نشان می‌دهد. کدام یک از ساختارهای زیر به احتمال زیاد در حین عمل جراحی رادیکال گردن آسیب دیده است؟"
$ws.Cells.Item(171, 1).Value = "در مورد کلسیم سرم همه موارد زیر صحیح است به جز"
$ws.Cells.Item(172, 1).Value = "طبق برنامه بهبودی تسریع‌شده در جراحی‌های روده بزرگ، کدام گزینه صحیح است؟"
$ws.Cells.Item(173, 1).Value = "لاتیسم به چه علتی رخ می‌دهد؟"
$ws.Cells.Item(174, 1).Value = "تومورهای تخمدان معمولاً از کدام منشأ می‌گیرند؟"
$ws.Cells.Item(175, 1).Value = "خط ارتعاشی خلفی در کجا قرار دارد؟"
$ws.Cells.Item(176, 1).Value = "کدام یک از موارد زیر در مورد iugr (تاخیر رشد داخل رحمی) صحیح است؟"
$ws.Cells.Item(177, 1).Value = "کدام یک از موارد زیر نوعی از شکستگی جمجمه نیست؟"
$ws.Cells.Item(178, 1).Value = "hcl توسط کدام سلول های معده ترشح می شود؟"
$ws.Cells.Item(179, 1).Value = "همه موارد زیر از علل فشار خون بالا همراه با هیپوکالمی هستند به جز،"
$ws.Cells.Item(180, 1).Value = "وزن طبیعی تیروئید با محتوای ید رژیم غذایی .................... تغییر می‌کند -"
$ws.Cells.Item(181, 1).Value = "شایع‌ترین محل کارسینومای مجرای ادرار در مردان کدام است؟"
$ws.Cells.Item(182, 1).Value = "علت پانکراتیت حاد همه موارد زیر به جز: مارس 2013 (g)"
$ws.Cells.Item(183, 1).Value = "کدام یک از ترشحات زیر دارای ph بسیار بالا است؟"
$ws.Cells.Item(184, 1).Value = "کاهش فشار خون باعث چه می‌شود؟"
$ws.Cells.Item(185, 1).Value = "یک کودک شش ساله برای جراحی انتخابی اورولوژی تحت بیهوشی عمومی برنامه‌ریزی شده است. او اجازه دسترسی وریدی به متخصص بیهوشی را نمی‌دهد. بهترین عامل استنشاقی برای القای بیهوشی کدام است؟"
$ws.Cells.Item(186, 1).Value = "کدام یک از عبارات زیر در مورد باکتریوفاژ صحیح است؟"
$ws.Cells.Item(187, 1).Value = "یک بیمار مرد 40 ساله با طحال بسیار بزرگ، حساسیت شدید جناغ و تعداد گلبول های سفید 85000 در هر میلی متر مکعب با درصد بالایی از میلوسیت ها و متامیلوسیت ها در بیمارستان بستری شده است. کدام یک از داروهای زیر برای بیماری او بهترین گزینه است؟"
$ws.Cells.Item(188, 1).Value = "کمپلکس اولیه در کدام یک از محل های زیر نشان دهنده سل مادرزادی است؟"
$ws.Cells.Item(189, 1).Value = "کدام یک از محیط‌های نگهداری اندام زیر دارای ph و اسمولاریته تقریباً مشابه با مقادیر ایده‌آل توصیه شده است؟"
$ws.Cells.Item(190, 1).Value = "تشکیل استخوان تقویتی (buttressing bone formation) یک"
$ws.Cells.Item(191, 1).Value = "کدام آمینو اسید دارای زنجیره جانبی بازی است؟"
$ws.Cells.Item(192, 1).Value = "کدام تومور استخوانی در اپی فیز رخ می‌دهد؟"
$ws.Cells.Item(193, 1).Value = "ناهنجاری های متابولیکی که در آن زردی دیده می شود، همه موارد زیر هستند به جز؟"
$ws.Cells.Item(194, 1).Value = "معایب روش مارشال مارکتی کرانتز در مقایسه با سایر گزینه های جراحی برای درمان بی اختیاری ادراری استرسی شامل کدام مورد است؟"
$ws.Cells.Item(195, 1).Value = "همه موارد زیر را می‌توان با ارکیوپکسی در کریپتورکیدیسم پیشگیری کرد به جز -"
$ws.Cells.Item(196, 1).Value = "میزان جریان بزاق تحریک‌نشده از غده پاروتید چقدر است؟"
$ws.Cells.Item(197, 1).Value = "هوردئولوم داخلی یک التهاب چرکی در کدام غده است؟"
$ws.Cells.Item(198, 1).Value = "یک بیمار جوان با آسیب شدید به دوازدهه پروگزیمال، سر پانکراس و مجرای صفراوی مشترک دیستال مراجعه می‌کند. روش انتخابی برای این بیمار باید چه باشد؟"
$ws.Cells.Item(199, 1).Value = "اولین نشانه افزایش فشار داخل جمجمه در معاینه ته چشم چیست؟"
$ws.Cells.Item(200, 1).Value = "آلودگی هوای داخل ساختمان در کدام مورد نقش ندارد؟"
$ws.Cells.Item(201, 1).Value = "کدام یک از موارد زیر با سندرم لوفلر مرتبط نیست؟"
$ws.Cells.Item(202, 1).Value = "علت اصلی مرگ و میر پریناتال در حالت بریچ فرانک ترم چیست؟"
$ws.Cells.Item(203, 1).Value = "نکروز انعقادی به چه علتی رخ می‌دهد؟"
$ws.Cells.Item(204, 1).Value = "بورلیا وینسنت چیست؟"
$ws.Cells.Item(205, 1).Value = "سطح بسیار بالای ترانس آمینازهای سرمی (alt/ast > 1000 iu/l) در موارد زیر دیده می‌شود به جز-"
$ws.Cells.Item(206, 1).Value = "سندرم وبر به دلیل ضایعات در کدام ناحیه ایجاد می‌شود؟"
$ws.Cells.Item(207, 1).Value = "بیشتر موارد فلج ایزوله عصب تروکلئار ناشی از کدام مورد است؟"
$ws.Cells.Item(208, 1).Value = "استین استراس چیست؟"
$ws.Cells.Item(209, 1).Value = "dna بسیار تکراری در کدام مورد مشاهده می‌شود؟"
$ws.Cells.Item(210, 1).Value = "تشخیص رادیوگرافی زیر چیست؟"
$ws.Cells.Item(211, 1).Value = "ایفوسفامید متعلق به کدام گروه از داروهای ضد سرطان است؟"
$ws.Cells.Item(212, 1).Value = "بی اختیاری مقعدی به دلیل درگیری کدام یک از موارد زیر است؟"
$ws.Cells.Item(213, 1).Value = "اریتم ندوزوم در همه موارد زیر دیده می‌شود به جز"
$ws.Cells.Item(214, 1).Value = "`"متاستاز استخوانی استئوسکلروتیک`" بیشتر در کدام کارسینوم دیده می‌شود:"
$ws.Cells.Item(215, 1).Value = "کدام یک از موارد زیر در مورد قلاب اینفرابالژ صحیح نیست؟"
$ws.Cells.Item(216, 1).Value = "سونوگرافی می‌تواند کیسه حاملگی را در چه زمانی زودتر تشخیص دهد؟"
$ws.Cells.Item(217, 1).Value = "در مورد نارکولپسی کدام گزینه نادرست است؟"
$ws.Cells.Item(218, 1).Value = "یک زن ۷۵ ساله در طول ۶ ماه گذشته درد مداوم ولی خفیفی در کمر، قفسه سینه راست، شانه چپ و ران چپ داشته است. اکنون به طور ناگهانی درد شدید و تیزی در ران چپش احساس می‌کند. در معاینه فیزیکی، درد شدیدی در لمس ران چپ مشاهده می‌شود و پای چپ کوتاه‌تر از پای راست است. رادیوگرافی پای چپ نشان‌دهنده شکستگی در ناحیه دیافیز فوقانی استخوان ران در یک ناحیه لیتیک ۵ سانتی‌متری است که تمام ضخامت استخوان را درگیر کرده است. اسکن استخوان چندین ناحیه با جذب افزایش‌یافته در استخوان ران چپ، لگن، مهره‌ها، دنده‌های سوم و چهارم راست، هومروس چپ فوقانی و اسکاپولای چپ نشان می‌دهد. آزمایش‌های آزمایشگاهی نشان‌دهنده کراتینین سرم ۰.۹ mg/dl، پروتئین کل ۶.۷ g/dl، آلبومین ۴.۵ g/dl، بیلی‌روبین کل ۱ mg/dl، ast 28 u/l، alt 22 u/l و آلکالین فسفاتاز ۲۰۲ u/l است. محتمل‌ترین تشخیص چیست؟"
$ws.Cells.Item(219, 1).Value = "کدام اسید آمینه صرفاً کتوژنیک است؟"
$ws.Cells.Item(220, 1).Value = "شایع‌ترین نوع دررفتگی مفصل ران کدام است؟"
$ws.Cells.Item(221, 1).Value = "بیمار با مشکلات عاطفی، افزایش ترشح بزاق، رنگ پریدگی مخاط دهان و تغییر رنگ خاکستری مایل به آبی لثه مراجعه می‌کند. این یافته‌ها بیشتر با کدام تشخیص بالینی سازگار است؟"
$ws.Cells.Item(222, 1).Value = "کدام گزینه در مورد مایکوباکتریوم‌های غیر از توبرکلوزیس صحیح است؟"
$ws.Cells.Item(223, 1).Value = "کدام یک بیشتر در ایجاد زگیل های تناسلی (فرج) نقش دارد؟"
$ws.Cells.Item(224, 1).Value = "علامت تریر اسکاتلندی در کدام نمای رادیوگرافی دیده می‌شود؟"
$ws.Cells.Item(225, 1).Value = "از میان گزینه‌های زیر، کدام یک از عوارض استروژن محسوب نمی‌شود؟"
$ws.Cells.Item(226, 1).Value = "یک فرد 33 ساله الکلی که تحت درمان att است، با افزایش آهن سرم و افزایش اشباع ترانسفرین مراجعه می‌کند. تشخیص؟"
$ws.Cells.Item(227, 1).Value = "تشخیص با ویژگی بالا چه چیزی را نشان می‌دهد؟"
$ws.Cells.Item(228, 1).Value = "dec (دی اتیل کاربامازین) برای درمان کدام بیماری استفاده می‌شود؟"
$ws.Cells.Item(229, 1).Value = "اختصاصی‌ترین علامت روان‌پریشی کدام است؟"
$ws.Cells.Item(230, 1).Value = "میزان حذف با سینتیک مرتبه اول متناسب است با -"
$ws.Cells.Item(231, 1).Value = "عامل ایجاد کننده کنژنکتیویت زاویه ای"
$ws.Cells.Item(232, 1).Value = "گچ آویزان برای چه موردی استفاده می‌شود؟"
$ws.Cells.Item(233, 1).Value = "یک پسر ۲ ساله با استریدور، تنگی نفس و مشکل در بلع به اورژانس مراجعه کرد. بیمار سابقه عفونت‌های مکرر قفسه سینه در گذشته را نیز داشت که به دلیل آن کودک به‌طور مکرر بستری می‌شد. اندوسکوپی دستگاه گوارش فوقانی انجام شد و طبیعی بود. یافته‌های آزمایشگاهی نیز طبیعی بودند. رادیوگرافی قفسه سینه انجام شد. چه تشخیصی محتمل‌تر است؟"
$ws.Cells.Item(234, 1).Value = "هیپوپلازی اندام‌ها و اسکار پوستی ناشی از کدام مورد است؟"
$ws.Cells.Item(235, 1).Value = "ناقل تب خونریزی دهنده کریمه کنگو کدام است؟"
$ws.Cells.Item(236, 1).Value = "یک بیمار زن 31 ساله از کاهش شنوایی دوطرفه در طول 5 سال گذشته شکایت دارد. در معاینه، پرده گوش طبیعی است و ادیوگرام کاهش شنوایی هدایتی دوطرفه را نشان می‌دهد. رفلکس‌های اکوستیک وجود ندارند. همه موارد زیر جزئی از درمان هستند به جز:"
$ws.Cells.Item(237, 1).Value = "بینایی در نوزادان چگونه ارزیابی می‌شود؟"
$ws.Cells.Item(238, 1).Value = "شاخص کورپوبازال برای تعیین کدام مورد مفید است؟"
$ws.Cells.Item(239, 1).Value = "کدام یک از موارد زیر به صورت یک توده در گوش که در صورت لمس به شدت خونریزی می‌کند، ظاهر می‌شود؟"
$ws.Cells.Item(240, 1).Value = "یک بیمار 31 ساله مبتلا به اسکیزوفرنی که به مدت 7 سال داروهای ضد روان‌پریشی مصرف می‌کرده است، برای معاینه چشم‌پزشکی مراجعه می‌کند. معاینه نشان‌دهنده وجود رنگدانه‌های شبکیه است. کدام یک از داروهای ضد روان‌پریشی زیر به احتمال زیاد توسط بیمار مصرف می‌شود؟"
$ws.Cells.Item(241, 1).Value = "در مورد اینترفرون‌ها کدام گزینه صحیح است؟"
$ws.Cells.Item(242, 1).Value = "کدام یک از موارد زیر محتوای بخش قدامی سوراخ ژوگولار است؟"
$ws.Cells.Item(243, 1).Value = "کدام یک از آزمایش‌های زیر در تشخیص نقص لوله عصبی کمک‌کننده است؟"
$ws.Cells.Item(244, 1).Value = "معنای کلمه اسکیزوفرنی چیست؟"
$ws.Cells.Item(245, 1).Value = "سلول‌های rs در کدام بیماری دیده می‌شوند؟"
$ws.Cells.Item(246, 1).Value = "یک خانم ۲۲ ساله فوت کرده است. کدام یک از یافته‌های پس از مرگ نشان‌دهنده این نیست که او زایمان کرده است؟"
$ws.Cells.Item(247, 1).Value = "بهترین عبارت در مورد مراکز تنفسی را انتخاب کنید:"
$ws.Cells.Item(248, 1).Value = "یک کودک 12 ساله با تب و لنفادنوپاتی گردنی مراجعه می‌کند. معاینه دهان یک غشای خاکستری روی لوزه راست که تا پایه قدامی گسترش یافته است را نشان می‌دهد. کدام یک از محیط‌های زیر برای کشت سواب گلو جهت شناسایی سریع پاتوژن ایده‌آل است؟"
$ws.Cells.Item(249, 1).Value = "کدام عامل ضد افسردگی به احتمال زیاد باعث ایجاد نعوظ پایدار (پریاپیسم) در یک بیمار مرد 40 ساله می‌شود؟"
$ws.Cells.Item(250, 1).Value = "آنتی‌ژن‌های abo در کدام یک از موارد زیر دیده نمی‌شوند؟"
$ws.Cells.Item(251, 1).Value = "افزایش ضربان قلب دقیقاً قبل از شروع ورزش به دلیل چیست؟"
$ws.Cells.Item(252, 1).Value = "ناحیه پیش حرکتی مربوط به کدام بخش است؟"
$ws.Cells.Item(253, 1).Value = "تمامی عبارات زیر در مورد سندرم داون صحیح هستند به جز؟"
$ws.Cells.Item(254, 1).Value = "شریان آپاندیکولار شاخه ای از کدام یک است؟"
$ws.Cells.Item(255, 1).Value = "کارسینوم فولیکولار تیروئید به دلیل جهش در کدام ژن ایجاد می‌شود؟"
$ws.Cells.Item(256, 1).Value = "حداقل پنوموتوراکس در کدام یک از نمای‌های زیر (رادیوگرافی قفسه سینه) بهتر دیده می‌شود؟"
$ws.Cells.Item(257, 1).Value = "اگر فیبرهای پورکینژ، که در قسمت دیستال اتصال دهلیزی-بطنی قرار دارند، به عنوان ضربانساز قلب عمل کنند، میزان ضربان قلب مورد انتظار چقدر خواهد بود؟"
$ws.Cells.Item(258, 1).Value = "کدام یک از موارد زیر به عنوان عامل خطر شایع برای بیماری عروق کرونر قلب محسوب نمی‌شود؟"
$ws.Cells.Item(259, 1).Value = "داروی انتخابی برای لوسمی میلوئید مزمن:"
$ws.Cells.Item(260, 1).Value = "خروج پلاسما و گلبول‌های سفید در التهاب حاد از کدام رگ‌ها صورت می‌گیرد؟"
$ws.Cells.Item(261, 1).Value = "اگر هم آب و هم غذا به طور کامل از فردی دریغ شود، مرگ معمولاً در چه زمانی رخ می‌دهد؟"
$ws.Cells.Item(262, 1).Value = "بیمار با وریدهای واریسی ضرباندار در اندام تحتانی مراجعه کرده است. محتمل‌ترین تشخیص چیست؟"
$ws.Cells.Item(263, 1).Value = "ذرات بی اثر با آنتی ژن یا آنتی بادی حساس می‌شوند. کدام یک از آزمایش‌های زیر به طور گسترده برای تشخیص سریع آنتی ژن‌های میکروبی (پنج دقیقه یا کمتر) استفاده می‌شود؟"
$ws.Cells.Item(264, 1).Value = "کدام یک از روش‌های زیر بهترین روش برای ارزیابی میزان مایعات دریافتی در بیمار پلی‌تراما است؟"
$ws.Cells.Item(265, 1).Value = "یک مواجهه بزرگ با پالپ در عرض 24 ساعت با چه روشی درمان می‌شود؟"
$ws.Cells.Item(266, 1).Value = "عفونت نهفته در موارد زیر رخ می دهد به جز-"
$ws.Cells.Item(267, 1).Value = "کدام گزینه در مورد کارسینوئیدهای برونشی صحیح است؟"
$ws.Cells.Item(268, 1).Value = "در یک شهر با جمعیت 10000 نفر، نسبت مردان به زنان 1 به 1 است. نسبت جنسیتی چیست؟"
$ws.Cells.Item(269, 1).Value = "تومور وارثین در کدام محل یافت می‌شود؟"
$ws.Cells.Item(270, 1).Value = "این علامت در کدام یک از موارد زیر مشاهده می‌شود؟"
$ws.Cells.Item(271, 1).Value = "تبدیل پره کالیکرئین به کالیکرئین به کدام فاکتور انعقادی نیاز دارد؟"
$ws.Cells.Item(272, 1).Value = "تشخیص آگاماگلوبولینمی وابسته به x در چه شرایطی باید مورد شک قرار گیرد؟"
$ws.Cells.Item(273, 1).Value = "مکانیسم‌های انتقال دارو شامل چه مواردی است؟"
$ws.Cells.Item(274, 1).Value = "محاسبه کمبود آب برای مرد 50 کیلوگرمی با سدیم=160 میلی‌اکیوالان در لیتر"
$ws.Cells.Item(275, 1).Value = "لنفانژیوگرافی پا چگونه انجام می‌شود؟"
$ws.Cells.Item(276, 1).Value = "کدام یک از موارد زیر به متابولیت فعال تبدیل نمی‌شود؟"
$ws.Cells.Item(277, 1).Value = "کدام یک از روش‌های جراحی زیر به عنوان زخم تمیز-آلوده در نظر گرفته می‌شود؟"
$ws.Cells.Item(278, 1).Value = "رفلکس مورو توسط چه چیزی آغاز می‌شود؟"
$ws.Cells.Item(279, 1).Value = "شایع‌ترین زیرگونه سرطان تیروئید کدام است؟"
$ws.Cells.Item(280, 1).Value = "کدام یک از موارد زیر در درمان کارسینومای تیروئید با تمایز خوب استفاده می‌شود؟"
$ws.Cells.Item(281, 1).Value = "آنتاگونیسم بین استیل کولین و آتروپین"
$ws.Cells.Item(282, 1).Value = "عارضه اصلی سیستوگاستروستومی برای کیست کاذب پانکراس چیست؟"
$ws.Cells.Item(283, 1).Value = "آزمایش پنی برای تشخیص چه چیزی انجام می‌شود؟"
$ws.Cells.Item(284, 1).Value = "تنظیم کننده سرعت تنفس (پیس میکر تنفسی):"
$ws.Cells.Item(285, 1).Value = "افزایش نسبت انسولین به گلوکاگون باعث چه چیزی می‌شود؟"
$ws.Cells.Item(286, 1).Value = "کدام یک از موارد زیر جفت غدد بزاقی نیست؟"
$ws.Cells.Item(287, 1).Value = "محتوای آلی رسوبات سوپراژینجیوال:"
$ws.Cells.Item(288, 1).Value = "کدام یک از علائم زیر نشان دهنده نارسایی احتقانی قلب (chf) در نارسایی است؟"
$ws.Cells.Item(289, 1).Value = "یک مرد 53 ساله برای معاینه بیمه به مطب پزشک مراجعه می‌کند. در طی گرفتن تاریخچه، شما متوجه می‌شوید که پدرش دیابت داشت و به دلیل نارسایی احتقانی قلب فوت کرده است. به بیمار گفته شده بود که این وضعیت به دلیل `"دیابت برنزی`" بوده است. بیمار فقط از خستگی مبهم و دردهای مفصلی شکایت دارد. شما نگران هستید که بیمار ممکن است در معرض خطر هموکروماتوز ارثی باشد. برای تأیید شک خود کدام آزمایش را درخواست می‌کنید؟"
$ws.Cells.Item(290, 1).Value = "بیمار با هیدروپس و اسکار قرنیه مراجعه کرده است. چه روش درمانی را برای او انتخاب می‌کنید؟"
$ws.Cells.Item(291, 1).Value = "تست مثبت رومبرگ با چشمان بسته، نقص در کدام قسمت را تشخیص می‌دهد؟"
$ws.Cells.Item(292, 1).Value = "آنافیلاکسی به چه چیزی اشاره دارد؟"
$ws.Cells.Item(293, 1).Value = "ترومباستنی گلانزمن چیست؟"
$ws.Cells.Item(294, 1).Value = "کدام ناحیه از مغز بیشتر دچار آتروفی در بیماری آلزایمر می‌شود؟"
$ws.Cells.Item(295, 1).Value = "از دست دادن ناگهانی و بدون درد بینایی در کدام یک از موارد زیر مشاهده می‌شود؟"
$ws.Cells.Item(296, 1).Value = "دسته کنت با کدام یک از موارد زیر مرتبط است؟"
$ws.Cells.Item(297, 1).Value = "کدام یک از انواع تومورهای تخمدان زیر بیشتر باعث بروز عوارض در دوران بارداری می‌شود؟"
$ws.Cells.Item(298, 1).Value = "بوی ادرار گربه نر در کدام یک از موارد زیر مشاهده می‌شود؟"
$ws.Cells.Item(299, 1).Value = "همه به جز یکی بر حداقل غلظت آلوئولی (mac) تأثیر می‌گذارند"
$ws.Cells.Item(300, 1).Value = "کدام یک از لایه‌های زیر جزو لایه‌های شبکیه نیست؟"
$ws.Cells.Item(301, 1).Value = "یک زن 36 ساله از ترشح خونی از نوک پستان به مدت 3 ماه شکایت دارد. در معاینه، یک ندول کوچک در عمق آرئول یافت می‌شود. لمس دقیق ناحیه نوک پستان-آرئول باعث ظهور خون در موقعیت ساعت 3 می‌شود. یافته‌های ماموگرافی طبیعی است. محتمل‌ترین تشخیص چیست؟"
$ws.Cells.Item(302, 1).Value = "گشاد شدن شکاف بین کندیلی استخوان ران در کدام مورد مشاهده می‌شود؟"
$ws.Cells.Item(303, 1).Value = "پدیده `"سرقت کرونری`" توسط کدام عامل ایجاد می‌شود؟"
$ws.Cells.Item(304, 1).Value = "کدام سندرم با دندان های اضافی همراه است؟"
$ws.Cells.Item(305, 1).Value = "کدام یک از موارد زیر از ویژگی‌های فیزیولوژیک عضله قلب نیست؟"
$ws.Cells.Item(306, 1).Value = "یک پسر یک ساله با هپاتوسپلنومگالی و تاخیر در مراحل رشدی مراجعه کرده است. بیوپسی کبد و مغز استخوان وجود هیستیوسیت‌هایی با مواد pas-مثبت مقاوم به دیاستاز در سیتوپلاسم را نشان داد. بررسی میکروسکوپ الکترونی این هیستیوسیت‌ها به احتمال زیاد وجود چه چیزی را نشان می‌دهد؟"
$ws.Cells.Item(307, 1).Value = "همه موارد زیر می‌توانند باعث انسداد سیستم هپاتوبیلیاری شوند به جز"
$ws.Cells.Item(308, 1).Value = "در منحنی تفکیک اکسی هموگلوبین، اشباع اکسیژن 95 درصد با __________ فشار اکسیژن در خون شریانی مطابقت دارد."
$ws.Cells.Item(309, 1).Value = "جریان خون رحم در ترم"
$ws.Cells.Item(310, 1).Value = "کدام گزینه در مورد آمیلوئیدوز صحیح نیست (غیر مرتبط)؟"
$ws.Cells.Item(311, 1).Value = "کدام یک از روش های جراحی زیر در درمان میوم همراه با خونریزی شدید قاعدگی کمترین اثر را دارد؟"
$ws.Cells.Item(312, 1).Value = "کیست زیرملتحمه در کدام یک از موارد زیر دیده می‌شود؟"
$ws.Cells.Item(313, 1).Value = "درمان اصلی یوئیت چیست؟"
$ws.Cells.Item(314, 1).Value = "ویروس هپاتیت c جزو کدام دسته است؟"
$ws.Cells.Item(315, 1).Value = "کدام گزینه در مورد مانیتورینگ پی‌اچ سرپایی نادرست است؟"
$ws.Cells.Item(316, 1).Value = "یک فرد 40 ساله الکلی دچار سرفه و تب می‌شود. عکس قفسه سینه سطح هوا-مایع را در بخش فوقانی لوب تحتانی راست نشان می‌دهد. محتمل‌ترین عامل اتیولوژیک کدام است؟"
$ws.Cells.Item(317, 1).Value = "همه موارد زیر باعث التهاب ملتحمه حاد می‌شوند به جز"
$ws.Cells.Item(318, 1).Value = "کدام یک از دندان های دائمی زیر، کمترین احتمال را برای داشتن دو ریشه دارد؟"
$ws.Cells.Item(319, 1).Value = "تمام موارد زیر از ویژگی های میخ کونچر (k) هستند که توسط یک جراح آلمانی برای فیکساسیون داخلی شکستگی های فمور طراحی شده است، به جز"
$ws.Cells.Item(320, 1).Value = "کدام یک از لیگامان‌های زیر از هایپراکستنشن (بازشدن بیش از حد) مفصل ران جلوگیری می‌کند؟"
$ws.Cells.Item(321, 1).Value = "کدام گزینه در مورد ویروس زیکا صحیح نیست؟"
$ws.Cells.Item(322, 1).Value = "شایع‌ترین محل طحال فرعی کجاست؟"
$ws.Cells.Item(323, 1).Value = "یک حساس کننده کانال کلسیمی که برای درمان نارسایی احتقانی قلب تأیید شده است کدام است؟"
$ws.Cells.Item(324, 1).Value = "در یک بیمار به شدت بیمار، افزودن اسیدهای آمینه به رژیم غذایی منجر به تعادل مثبت نیتروژن می‌شود. مکانیسم این اثر چیست؟"
$ws.Cells.Item(325, 1).Value = "فرد 60 ساله ای با سابقه آنژین و تنگی نفس در طول هفته گذشته مراجعه می کند. خون گرفته شده رنگ قهوه ای مایل به قرمز غلیظ دارد. تشخیص چیست؟"
$ws.Cells.Item(326, 1).Value = "در مورد غدد پاراتیروئید صحیح است به جز"
$ws.Cells.Item(327, 1).Value = "کدام یک از موارد زیر علت دیلاسراسیون‌ها است؟"
$ws.Cells.Item(328, 1).Value = "در تصویر بالا، مرز خلفی مثلث تراگتمن توسط چه ساختاری تشکیل شده است؟"
$ws.Cells.Item(329, 1).Value = "نامطلوب‌ترین شکل ترمیم شکستگی ریشه چیست؟"
$ws.Cells.Item(330, 1).Value = "در dic کدام گزینه صحیح است:"
$ws.Cells.Item(331, 1).Value = "موج 'c' در jvp"
$ws.Cells.Item(332, 1).Value = "داروهای ضد رتروویروسی برای بیماران hiv/aids زمانی شروع می‌شوند که تعداد سلول‌های cd4 به چه میزان باشد؟"
$ws.Cells.Item(333, 1).Value = "در مورد واکنش آگلوتیناسیون همه موارد زیر صحیح هستند به جز: سپتامبر 2007"
$ws.Cells.Item(334, 1).Value = "دیورتیکول معده  
الف) درد علامت اصلی است  
ب) معمولاً در انتهای کاردیاک قرار دارد  
ج) معمولاً در سطح خلفی است  
د) اینورژن درمان رضایتبخش است  
ه) همه موارد فوق"
$ws.Cells.Item(335, 1).Value = "کدام یک از موارد زیر جزء بیماری‌های میلوپرولیفراتیو نیست؟"
$ws.Cells.Item(336, 1).Value = "هنگامی که از کف دست به عنوان مقیاس اندازه‌گیری برای تخمین مساحت سوختگی استفاده می‌شود، معادل چه درصدی از سطح بدن است؟"
$ws.Cells.Item(337, 1).Value = "کدام یک از موارد زیر شایع‌ترین علت تنگی انسدادی مجرای ادرار غشایی است؟"
$ws.Cells.Item(338, 1).Value = "کدام یک از موارد زیر مهارکننده فسفوریلاسیون اکسیداتیو نیست؟"
$ws.Cells.Item(339, 1).Value = "هیدروکسی اتیل نشاسته به عنوان چه چیزی استفاده می‌شود؟"
$ws.Cells.Item(340, 1).Value = "اجسام مالوری هیالین به طور مشخص در کدام بیماری دیده می‌شوند؟"
$ws.Cells.Item(341, 1).Value = "کدامیک از موارد زیر باعث کلسیفیکاسیون پریکارد نمی‌شود؟"
$ws.Cells.Item(342, 1).Value = "در روش نمونه‌گیری خوشه‌ای epi توصیه شده توسط who برای ارزیابی پوشش ایمن‌سازی اولیه، محدوده سنی کودکان مورد بررسی کدام است؟"
$ws.Cells.Item(343, 1).Value = "در سن 35 سالگی، خطر مرتبط با سن مادر برای سندرم داون در بارداری تک‌قلویی در سه‌ماهه دوم چقدر است؟"
$ws.Cells.Item(344, 1).Value = "فاز پلاتوی عضله بطنی به دلیل باز شدن کدام کانال است؟"
$ws.Cells.Item(345, 1).Value = "روش جلوگیری از بارداری که در دوران شیردهی منع مصرف دارد کدام است؟"
$ws.Cells.Item(346, 1).Value = "یک خانم ۲۶ ساله نخستزا در هفته ۳۰ بارداری در معاینه بالینی، مقادیر فشار خون ۱۴۲/۱۰۰ میلیمتر جیوه، ۱۵۰/۹۴ میلیمتر جیوه و ۱۵۰/۱۰۰ میلیمتر جیوه در فواصل ۶ ساعته دارد. اقدام بعدی در مدیریت این بیمار چیست؟"
$ws.Cells.Item(347, 1).Value = "مفهوم `"پزشکی اجتماعی شده`" اولین بار توسط کدام کشور ارائه شد؟"
$ws.Cells.Item(348, 1).Value = "اتونفرکتومی در کدام مورد مشاهده می‌شود؟"
$ws.Cells.Item(349, 1).Value = "یک مرد 25 ساله با تب، سرفه، خلط و تنگی نفس به مدت 2 ماه مراجعه کرده است. سی تی اسکن قفسه سینه با کنتراست، ضایعات فیبروتیک دوطرفه در لوب فوقانی و غدد لنفاوی نکروتیک بزرگ شده مدیاستین با تقویت حاشیه ای محیطی را نشان داد. کدام یک از موارد زیر محتمل ترین تشخیص است؟"
$ws.Cells.Item(350, 1).Value = "هیپوفسفاتمی ناشی از کدام مورد زیر است؟"
$ws.Cells.Item(351, 1).Value = "روش انتخابی برای بررسی زنان در معرض خطر بالای سرطان پستان چیست؟"
$ws.Cells.Item(352, 1).Value = "دانه‌های بربک در کدام مورد مشاهده می‌شود؟"
$ws.Cells.Item(353, 1).Value = "صدای نرم s1 در تمام موارد زیر دیده می‌شود به جز"
$ws.Cells.Item(354, 1).Value = "یک کودک با مشکل عدم رشد کافی به درمانگاه اطفال مراجعه می‌کند. پزشک به والدین در مورد بهبود تغذیه کودک مشاوره می‌دهد، که مادر به شدت تأکید می‌کند که آنها از طرف خود همه چیز را به ویژه تغذیه به طور کامل رعایت کرده‌اند، اما نتیجه‌ای نداشته است. سپس پزشک درخواست چند آزمایش می‌کند. تشخیص چیست؟"
$ws.Cells.Item(355, 1).Value = "جهش uac به uag"
$ws.Cells.Item(356, 1).Value = "روش تشخیصی انتخابی برای دیورتیکول زنکر چیست؟"
$ws.Cells.Item(357, 1).Value = "مهارکننده های تیروزین کیناز به عنوان درمان خط اول در کدام مورد استفاده می‌شوند؟"
$ws.Cells.Item(358, 1).Value = "vitremer چیست؟"
$ws.Cells.Item(359, 1).Value = "کدام یک از بی‌حس‌کننده‌های موضعی زیر یک انقباض‌دهنده عروق است؟"
$ws.Cells.Item(360, 1).Value = "بیماری جانسن چیست؟"
$ws.Cells.Item(361, 1).Value = "یک زن ۷۰ ساله دو روز پیش به دلیل شکستگی لگن پس از زمین خوردن، تحت عمل جراحی لگن قرار گرفته است. هیچ سابقه‌ای از جراحی یا مصرف دارو در گذشته وجود ندارد. در طول ۲۴ ساعت گذشته، او از ناراحتی و اتساع شکم شکایت دارد. در معاینه: بدون تب، فشار خون ۱۴۰/۸۰ میلی‌متر جیوه، ضربان قلب ۱۱۰ در دقیقه، تعداد تنفس ۱۶ در دقیقه. شکم او متسع و طبل‌گونه با عدم وجود صداهای روده است. هیچ حساسیت برگشتی وجود ندارد. تصویر رادیوگرافی شکم در حالت ایستاده نشان داده شده است: محتمل‌ترین تشخیص در این مورد چیست؟"
$ws.Cells.Item(362, 1).Value = "استخوان نخودی با کدام استخوان مفصل می‌شود؟"
$ws.Cells.Item(363, 1).Value = "ناحیه خودمختار عصب اولنار"
$ws.Cells.Item(364, 1).Value = "ثبت تولد باید در چه مدت زمانی انجام شود؟"
$ws.Cells.Item(365, 1).Value = "کدام یک از موارد زیر یک اختلال مغلوب وابسته به اکس نیست؟"
$ws.Cells.Item(366, 1).Value = "فعالیت کدام یک از آنزیم های زیر در حالت گرسنگی افزایش نمی یابد؟"
$ws.Cells.Item(367, 1).Value = "متداول‌ترین روش برای بلوک شبکه بازویی چیست؟"
$ws.Cells.Item(368, 1).Value = "در طول روزه‌داری، یک فاز کوتاه از انقباضات شدید و متوالی در معده شروع می‌شود و به تدریج به ایلئوم مهاجرت می‌کند. ترشح کدام یک از هورمون‌های روده‌ای زیر به احتمال زیاد مسئول این اثر مشاهده‌شده است؟"
$ws.Cells.Item(369, 1).Value = "یک کودک دختر با خونریزی واژینال به مدت 4 روز پس از تولد مراجعه می‌کند، چه اقداماتی باید انجام شود؟"
$ws.Cells.Item(370, 1).Value = "کدام گزینه درباره فوسفنی توئین نادرست است؟"
$ws.Cells.Item(371, 1).Value = "در مخزن مالاریا، انگل به چه شکلی باقی می‌ماند؟"
$ws.Cells.Item(372, 1).Value = "hutch دیورتیکول در کدام قسمت دیده می‌شود؟"
$ws.Cells.Item(373, 1).Value = "یک زن که 3 سال پیش جراحی سرطان سینه چپ را انجام داده است، اکنون دچار گره آبی رنگ در همان سمت شده است."
$ws.Cells.Item(374, 1).Value = "یک مرد ۷۰ ساله که در طول ۵۰ سال گذشته تنباکو می‌جویده، با سابقه شش ماهه یک ضایعه بزرگ، قارچ‌مانند، نرم و پاپیلاری در حفره دهان مراجعه می‌کند. این ضایعه به داخل استخوان فک پایین نفوذ کرده است. غدد لنفاوی قابل لمس نیستند. دو نمونه بیوپسی گرفته شده از خود ضایعه، پاپیلوماتوز خوش‌خیم با هایپرکراتوز و آکانتوز را نشان می‌دهد که به بافت‌های زیرین نفوذ کرده است. محتمل‌ترین تشخیص چیست؟"
$ws.Cells.Item(375, 1).Value = "عضلات جویدن از کدام قوس مشتق می‌شوند؟"
$ws.Cells.Item(376, 1).Value = "بتا-لاکتامازها چگونه باعث مقاومت به پنیسیلین‌ها و سفالوسپورین‌ها می‌شوند؟"
$ws.Cells.Item(377, 1).Value = "کدام یک از موارد زیر در مورد عفونت بیمارستانی صحیح نیست؟"
$ws.Cells.Item(378, 1).Value = "گره های لنفاوی سابمنتال توسط همه موارد زیر تخلیه می شوند به جز-"
$ws.Cells.Item(379, 1).Value = "یک بیمار پس از انفارکتوس میوکارد دچار نارسایی میترال می‌شود. علت احتمالی این عارضه چیست؟"
$ws.Cells.Item(380, 1).Value = "معمولاً کاسپ پنجم در کدام یک دیده نمی‌شود؟"
$ws.Cells.Item(381, 1).Value = "مینوسایکلین باعث تغییر رنگ کدام قسمت می‌شود؟"
$ws.Cells.Item(382, 1).Value = "بیماری جانسن چیست؟"
$ws.Cells.Item(383, 1).Value = "کدام یک از گزینه‌های زیر توسط mhc ii ارائه نمی‌شود؟"
$ws.Cells.Item(384, 1).Value = "تمام موارد زیر از camp به عنوان پیام‌رسان دوم در عملکرد فیزیولوژیک خود استفاده می‌کنند، به جز:"
$ws.Cells.Item(385, 1).Value = "در مورد gist همه موارد زیر صحیح است به جز -"
$ws.Cells.Item(386, 1).Value = "تغییر رنگ نارنجی قرنیه ناشی از چیست؟"
$ws.Cells.Item(387, 1).Value = "برای کدام بدخیمی، پرتودرمانی با شدت تعدیل‌شده (imrt) مناسب‌ترین گزینه است؟"
$ws.Cells.Item(388, 1).Value = "تیفوس اپیدمیک ناشی از کدام عامل است؟"
$ws.Cells.Item(389, 1).Value = "کدام گزینه در مورد بیماری ویلسون صحیح است؟"
$ws.Cells.Item(390, 1).Value = "همه موارد زیر در مورد دیورتیکول مکل صحیح هستند به جز"
$ws.Cells.Item(391, 1).Value = "در مورد آسیب عصبی محیطی در اندام فوقانی کدام گزینه صحیح است؟"
$ws.Cells.Item(392, 1).Value = "آنتی بادی در لوپوس ناشی از دارو"
$ws.Cells.Item(393, 1).Value = "سلول t در رد پیوند کدام آنتی‌ژن را تشخیص می‌دهد؟"
$ws.Cells.Item(394, 1).Value = "مکانیسم اثر d-tubocurarine چیست؟"
$ws.Cells.Item(395, 1).Value = "پروتئین‌های بنس جونز در همه موارد زیر یافت می‌شوند به جز"
$ws.Cells.Item(396, 1).Value = "کدام یک از موارد زیر فیلتراسیون در انتهای شریانی بستر مویرگی را تسهیل می‌کند؟"
$ws.Cells.Item(397, 1).Value = "کدام یک از موارد زیر در مورد مایعات بدن نادرست است؟"
$ws.Cells.Item(398, 1).Value = "تکنیکی که دو ماده مایع (یا خمیری) را بدون همزدن مکانیکی به مخلوطی همگن تبدیل می‌کند، به چه نام شناخته می‌شود؟"
$ws.Cells.Item(399, 1).Value = "رژیم مورد استفاده برای مدیریت انتظاری جفت سرراهی چیست؟"
$ws.Cells.Item(400, 1).Value = "کدام یک از داروهای ضد صرع زیر را می‌توان برای درمان زردی نوزادان استفاده کرد؟"
$ws.Cells.Item(401, 1).Value = "گرانولومای اینگوئیناله توسط کدام عامل ایجاد می‌شود؟"
$ws.Cells.Item(402, 1).Value = "طبق دستورالعمل‌های rntcp اولین اقدام در مورد مشکوک به سل چیست؟"
$ws.Cells.Item(403, 1).Value = "تمام موارد زیر می‌توانند برای کند کردن فرآیند سفت شدن خمیر قالب‌گیری اکسید روی اوژنول استفاده شوند به جز"
$ws.Cells.Item(404, 1).Value = "اگر یک کودک 11 ماهه دو دوز واکسن dpi و فلج اطفال را دریافت کرده باشد و پس از 5 ماه از آخرین دوز برای ایمن‌سازی بیشتر مراجعه کند، چه باید کرد؟"
$ws.Cells.Item(405, 1).Value = "بهترین توصیف برای سینوزوئیدها کدام است؟"
$ws.Cells.Item(406, 1).Value = "در طول احیای قلبی ریوی پیشرفته (cpcr) در acls، تغییر زیر در موج کاپنوگرافی مشاهده می‌شود. این نشان‌دهنده چیست؟"
$ws.Cells.Item(407, 1).Value = "یک کودک مبتلا به سوءتغذیه از وضعیت اجتماعی-اقتصادی ضعیف که در مناطق شلوغ و کثیف زندگی می‌کند، با یک ندول در اطراف لیمبوس همراه با هیپرمی ملتحمه اطراف در چشم چپ خود مراجعه می‌کند. همچنین تورم غدد لنفاوی زیر بغل و گردن در او مشاهده می‌شود. کدام یک از موارد زیر محتمل‌ترین تشخیص است؟"
$ws.Cells.Item(408, 1).Value = "یک زن باردار برای سومین بار با شکایت از درد ناگهانی در قسمت پایین شکم سمت راست پس از ۵ هفته آمنوره (قطع قاعدگی) بستری شده است. ضربان نبض او ۱۳۰ و فشار خون ۸۰/۵۰ mmhg است. سونوگرافی واژینال نشان دهنده مقدار زیادی مایع آزاد در لگن و حفره رحمی خالی است. مرحله بعدی چیست؟"
$ws.Cells.Item(409, 1).Value = "تیموما با همه موارد زیر مرتبط است، به جز:"
$ws.Cells.Item(410, 1).Value = "شمشیر دو لبه کدام است؟"
$ws.Cells.Item(411, 1).Value = "شاخص کیفیت خدمات سلامت مادر و کودک (mch) –"
$ws.Cells.Item(412, 1).Value = "وزن جنین سقط شده چقدر است؟"
$ws.Cells.Item(413, 1).Value = "در مورد پدیده لوسیو کدام گزینه صحیح است؟"
$ws.Cells.Item(414, 1).Value = "علامت هگار در چه زمانی زودتر دیده می‌شود؟"
$ws.Cells.Item(415, 1).Value = "اندامک‌هایی که نقش محوری در آپوپتوز دارند کدامند؟"
$ws.Cells.Item(416, 1).Value = "سهم اصلی در گلوکونئوژنز مربوط به کدام مورد است؟"
$ws.Cells.Item(417, 1).Value = "شایع‌ترین محل متاستاز سرطان دهانه رحم کدام است؟"
$ws.Cells.Item(418, 1).Value = "کدام گزینه در مورد پلاسمیدهای باکتریایی نادرست است؟"
$ws.Cells.Item(419, 1).Value = "سندرم ولن نشان دهنده کدام یک از موارد زیر است؟"
$ws.Cells.Item(420, 1).Value = "علامت لیسکر در کدام مورد مشاهده می‌شود؟"
$ws.Cells.Item(421, 1).Value = "اولین سلولی که به دلیل کموتاکسی به سمت زخم مهاجرت می‌کند تا فرآیند بهبود زخم را آغاز کند کدام است؟"
$ws.Cells.Item(422, 1).Value = "همه موارد زیر از دلایل درد در فیبروم زیرمخاطی هستند، به جز کدام یک؟"
$ws.Cells.Item(423, 1).Value = "کدام یک از موارد زیر بهترین بیمار برای دندان مصنوعی است؟"
$ws.Cells.Item(424, 1).Value = "همه موارد زیر از علل باز شدن زخم پس از عمل هستند به جز-"
$ws.Cells.Item(425, 1).Value = "کم خونی همراه با رتیکولوسیتوز در کدام مورد مشاهده می‌شود؟"
$ws.Cells.Item(426, 1).Value = "کدام یک از شاخه‌های قسمت گردنی شریان فاسیال نیست؟"
$ws.Cells.Item(427, 1).Value = "کدام یک از گزینه‌های زیر زیباترین پونتیک است؟"
$ws.Cells.Item(428, 1).Value = "سندرمی که با هایپرکراتوز دست و پا همراه با تخریب شدید بافت های پریودنتال در دندان های شیری و دائمی و کلسیفیکاسیون دورا مشخص می شود."
$ws.Cells.Item(429, 1).Value = "سطح پیشنهادی سر و صدایی که افراد می‌توانند بدون آسیب به شنوایی تحمل کنند، چیست؟"
$ws.Cells.Item(430, 1).Value = "یک آزمایش ساده باکتریایی برای تشخیص مواد سرطان‌زای جهش‌زا چیست؟"
$ws.Cells.Item(431, 1).Value = "عفونت‌های منتقل شده به نوزاد در هنگام زایمان:"
$ws.Cells.Item(432, 1).Value = "سه‌گانه پیگمانتاسیون پری‌اورال، سرطان روده و سرطان پستان در کدام مورد مشاهده می‌شود؟"
$ws.Cells.Item(433, 1).Value = "تمامی موارد زیر با سوفل سیستولیک اولیه مرتبط هستند، به جز:"
$ws.Cells.Item(434, 1).Value = "affiliation به چه معناست؟"
$ws.Cells.Item(435, 1).Value = "برای درمان ادم مغزی ناشی از تومور مغزی، کورتیکواستروئید ترجیحی کدام است؟"
$ws.Cells.Item(436, 1).Value = "کدام یک از دستگاه‌های زیر، یک دستگاه تحویل کنترل‌شده است که برای تحویل غلظت ثابتی از اکسیژن استفاده می‌شود؟"
$ws.Cells.Item(437, 1).Value = "جمعیتی 1000 نفری که در آن شیوع بیماری xyz برابر 50٪ است. از تست غربالگری abc برای تشخیص بیماری استفاده شده که حساسیت و ویژگی آن هر دو 50٪ است. ارزش پیشبینیکننده مثبت تست غربالگری abc چقدر است؟"
$ws.Cells.Item(438, 1).Value = "اصطلاح کاتاتونیا توسط چه کسی استفاده شد؟"
$ws.Cells.Item(439, 1).Value = "توالی متفاوت اسیدهای آمینه با ساختار مشابه پروتئین‌ها نمونه‌ای از چیست؟"
$ws.Cells.Item(440, 1).Value = "در کدام یک از مسمومیت‌های زیر، قولنج شکمی یک ویژگی بارز است؟"
$ws.Cells.Item(441, 1).Value = "پروتئین اتصال دهنده هموگلوبین کدام است؟"
$ws.Cells.Item(442, 1).Value = "نادرست در مورد باسیلوس آنتراسیس کدام است؟"
$ws.Cells.Item(443, 1).Value = "آزوتمی پیش کلیوی با کدام یک از ویژگی های زیر مرتبط است؟"
$ws.Cells.Item(444, 1).Value = "آسیب به مجرای ادرار مردانه در زیر غشاء پرینه باعث تجمع ادرار در کدام ناحیه می‌شود؟"
$ws.Cells.Item(445, 1).Value = "جهش‌های مختلف در یک جایگاه ژنتیکی یکسان که باعث فنوتیپ مشابه یا یکسان می‌شوند، چیست؟"
$ws.Cells.Item(446, 1).Value = "کدام یک از موارد زیر پیشگیری اولیه محسوب می‌شود؟"
$ws.Cells.Item(447, 1).Value = "انتهای قدامی طحال توسط چه چیزی نگه داشته می‌شود؟"
$ws.Cells.Item(448, 1).Value = "آسیب به بخش شکمی-جانبی نخاع منجر به از دست دادن (در زیر سطح آسیب) کدام یک از موارد زیر می‌شود؟"
$ws.Cells.Item(449, 1).Value = "درد در گوش در چه سطحی رخ می‌دهد؟"
$ws.Cells.Item(450, 1).Value = "کدام مورد در سندرم ویلیام دیده نمی‌شود؟"
$ws.Cells.Item(451, 1).Value = "بار آهن در همه موارد زیر رخ می‌دهد، به جز -"
$ws.Cells.Item(452, 1).Value = "آگورافوبیا چیست؟"
$ws.Cells.Item(453, 1).Value = "اگر حرکت پیش‌جراحی قطعات فک بالا در یک نوزاد مبتلا به شکاف لب و کام نشان داده شود، از چه زمانی شروع می‌شود؟"
$ws.Cells.Item(454, 1).Value = "کدام یک از اسیدهای آمینه، هیدروفوبیک (آب‌گریز) است؟"
$ws.Cells.Item(455, 1).Value = "استخوان رکابی از کدام قوس تشکیل میشود؟"
$ws.Cells.Item(456, 1).Value = "کدام موقعیت، محل قرارگیری شبکه سلیاک را به بهترین شکل توصیف می‌کند؟"
$ws.Cells.Item(457, 1).Value = "تمام موارد زیر از موارد استفاده غیرپیشگیری از بارداری قرص های ocp هستند، به جز:"
$ws.Cells.Item(458, 1).Value = "یک نوزاد با تشخیص پنتالوژی فالوت ممکن است کدام یک از ضایعات زیر را داشته باشد؟"
$ws.Cells.Item(459, 1).Value = "همه موارد زیر در سوختگی‌های آنتموئم دیده می‌شوند، به جز:"
$ws.Cells.Item(460, 1).Value = "جمعیت تحت پوشش یک مرکز بهداشت اولیه (phc) در مناطق کوهستانی چقدر است؟"
$ws.Cells.Item(461, 1).Value = "نرخ بازسازی کام و زبان چقدر است؟"
$ws.Cells.Item(462, 1).Value = "کدام یک از موارد زیر در مورد مارپیچ آلفا صحیح نیست؟"
$ws.Cells.Item(463, 1).Value = "داروی مورد استفاده در پیشگیری از h. influenza کدام است؟"
$ws.Cells.Item(464, 1).Value = "اچ کردن مینا چند الگو ایجاد می‌کند؟"
$ws.Cells.Item(465, 1).Value = "فشار طبیعی انقباضات رحمی بین 190 تا 300 واحد است. در اینجا به کدام واحد اشاره شده است؟"
$ws.Cells.Item(466, 1).Value = "آگلوتینین گرم در همه موارد زیر یافت می‌شود به جز:"
$ws.Cells.Item(467, 1).Value = "یک دختر ۳ ساله پس از فرو کردن مداد به مجرای گوش خارجی خود، پرده گوشش پاره شد. او به صورت اورژانسی به بخش اورژانس منتقل شد. معاینه فیزیکی نشان دهنده درد در گوش و چند قطره خون در مجرای شنوایی خارجی بود. نگرانی این بود که ممکن است به عصبی که عمدتاً سطح خارجی پرده گوش را عصب‌دهی می‌کند، آسیب رسیده باشد. کدام یک از آزمایش‌های زیر به احتمال زیاد در معاینه فیزیکی برای بررسی آسیب به این عصب انجام می‌شود؟"
$ws.Cells.Item(468, 1).Value = "یک کودک با مننژیت باکتریایی عودکننده باید از نظر بالینی برای حضور چه موردی ارزیابی شود؟"
$ws.Cells.Item(469, 1).Value = "از مواد منبسط کننده پلاسما در چه مواردی استفاده می‌شود؟"
$ws.Cells.Item(470, 1).Value = "s100 یک مارکر مورد استفاده در تشخیص تمام موارد زیر است به جز"
$ws.Cells.Item(471, 1).Value = "لکه‌های میلیاری ریه در همه موارد زیر دیده می‌شود به جز"
$ws.Cells.Item(472, 1).Value = "ضخیم‌شدگی موضعی و منتشر دیواره کیسه صفرا همراه با بازتاب‌های با دامنه بالا و آرتیفکت‌های «دنباله‌دار» در سونوگرافی، تشخیص کدام بیماری را مطرح می‌کند؟"
$ws.Cells.Item(473, 1).Value = "کدام یک از موارد زیر تحت تأثیر کلرزنی معمولی قرار نمی‌گیرد؟"
$ws.Cells.Item(474, 1).Value = "در مورد pecoma (تومور سلول اپی تلیوئید پری واسکولار) کدام گزینه صحیح است؟"
$ws.Cells.Item(475, 1).Value = "سندرم تیتزه معمولاً در غضروف دنده ای کدام ناحیه ایجاد می شود؟"
$ws.Cells.Item(476, 1).Value = "کیسه زرد توسط چه چیزی از بین می‌رود؟"
$ws.Cells.Item(477, 1).Value = "پلاسمای بیمار مبتلا به پرکاری تیروئید حاوی مقدار اضافی از کدام مورد است؟"
$ws.Cells.Item(478, 1).Value = "ریشه عصبی عصب پوستی داخلی ران"
$ws.Cells.Item(479, 1).Value = "تمامی عبارات زیر در مورد آرتریت نقرسی اولیه صحیح هستند، به جز"
$ws.Cells.Item(480, 1).Value = "اپی فیز کششی چیست؟"
$ws.Cells.Item(481, 1).Value = "کاهش شنوایی 65 دسی‌بل، درجه ناشنوایی چیست؟"
